# Update "System Info" worksheet: append 5 new data rows (43-47) mirroring
# the existing data rows, with N column continuing the sequence and the
# Specification column reflecting updated memory usage readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (which currently ends at row 42).
$newRows = @(
    @{ N = 42; Spec = "Total Memory: 15.86 GB, Used Memory: 6.97 GB, Total Disk Space: 237.84 GB" },
    @{ N = 43; Spec = "Total Memory: 15.86 GB, Used Memory: 7.26 GB, Total Disk Space: 237.84 GB" },
    @{ N = 44; Spec = "Total Memory: 15.86 GB, Used Memory: 7.43 GB, Total Disk Space: 237.84 GB" },
    @{ N = 45; Spec = "Total Memory: 15.86 GB, Used Memory: 7.33 GB, Total Disk Space: 237.84 GB" },
    @{ N = 46; Spec = "Total Memory: 15.86 GB, Used Memory: 7.35 GB, Total Disk Space: 237.84 GB" }
)

$username = "Sovan.Souern"
$serial = "1L0N1W2"
$model = "AMD64"
$manufacturer = "Windows"
$assetTag = "PNCL114"
$remark = "AT/AT COMPATIBLE"

$lastRow = 42
$startRow = $lastRow + 1
$endRow = $startRow + $newRows.Count - 1

# Copy formatting from the last existing row down across the new rows so
# the appended rows look consistent with the rest of the table.
$ws.Range("A$lastRow`:H$lastRow").Copy() | Out-Null
$ws.Range("A$startRow`:H$endRow").PasteSpecial(-4122) | Out-Null

$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row.N
    $ws.Cells.Item($r, 2).Value2 = $row.Spec
    $ws.Cells.Item($r, 3).Value2 = $username
    $ws.Cells.Item($r, 4).Value2 = $serial
    $ws.Cells.Item($r, 5).Value2 = $model
    $ws.Cells.Item($r, 6).Value2 = $manufacturer
    $ws.Cells.Item($r, 7).Value2 = $assetTag
    $ws.Cells.Item($r, 8).Value2 = $remark
    $r++
}

$excel.CutCopyMode = 0

Write-Host "Added rows $startRow to $endRow"
